$wb = $excel.ActiveWorkbook

# --- Add the new "dec 17" sheet as the first tab ---------------------------
$new = $wb.Worksheets.Add()
$new.Name = "dec 17"
$new.Move($wb.Worksheets.Item(1))

# Header / date row
$new.Range("A2").Value = 43099
$new.Range("A2").NumberFormat = "d-mmm-yy"
$new.Range("B2").Value = "Score"
$new.Range("C2").Value = "Fairway"
$new.Range("D2").Value = "GIR"
$new.Range("E2").Value = "Putts"
$new.Range("F2").Value = "Comment"

# Hole-by-hole data
$new.Range("A3").Value = "Hole 1"
$new.Range("B3").Value = 4
$new.Range("C3").Value = "S"
$new.Range("E3").Value = 1

$new.Range("A4").Value = "Hole 2"
$new.Range("B4").Value = 3
$new.Range("E4").Value = 1

$new.Range("A5").Value = "Hole 3"
$new.Range("B5").Value = 5
$new.Range("C5").Value = "R"
$new.Range("E5").Value = 2

$new.Range("A6").Value = "Hole 4"
$new.Range("B6").Value = 5
$new.Range("C6").Value = "R"
$new.Range("E6").Value = 1

$new.Range("A7").Value = "Hole 5"
$new.Range("B7").Value = 3
$new.Range("E7").Value = 2

$new.Range("A8").Value = "Hole 6"
$new.Range("B8").Value = 5
$new.Range("C8").Value = "S"
$new.Range("E8").Value = 2

$new.Range("A9").Value = "Hole 7"
$new.Range("B9").Value = 5
$new.Range("C9").Value = "S"
$new.Range("E9").Value = 1

$new.Range("A10").Value = "Hole 8"
$new.Range("B10").Value = 4
$new.Range("E10").Value = 2

$new.Range("A11").Value = "Hole 9"
$new.Range("B11").Value = 5
$new.Range("C11").Value = "S"
$new.Range("E11").Value = 2

$new.Range("A12").Value = "Hole 10"
$new.Range("B12").Value = 4
$new.Range("C12").Value = "L"
$new.Range("E12").Value = 2

$new.Range("A13").Value = "Hole 11"
$new.Range("B13").Value = 5
$new.Range("E13").Value = 3

$new.Range("A14").Value = "Hole 12"
$new.Range("B14").Value = 5
$new.Range("C14").Value = "R"
$new.Range("E14").Value = 2

$new.Range("A15").Value = "Hole 13"
$new.Range("A16").Value = "Hole 14"

$new.Range("A17").Value = "Hole 15"
$new.Range("B17").Value = 3
$new.Range("E17").Value = 2

$new.Range("A18").Value = "Hole 16"
$new.Range("A19").Value = "Hole 17"

$new.Range("A20").Value = "Hole 18"
$new.Range("B20").Value = 5
$new.Range("C20").Value = "R"
$new.Range("E20").Value = 2

# Totals row
$new.Range("B21").Formula = "=SUM(B3:B20)"
$new.Range("E21").Formula = "=SUM(E3:E20)"

# Reserved (blank, date-formatted) rows for the next round not yet logged
$new.Range("A23").NumberFormat = "d-mmm-yy"
$new.Range("A24").NumberFormat = "d-mmm-yy"

$new.Columns.Item(1).ColumnWidth = 10.109375

# sheet view: dec17 is the tab on display, scrolled down, B32 selected
$new.Application.ActiveWindow.ScrollRow = 13
[void]$new.Range("B32").Select()

# --- "jan 18" tab (formerly the first tab) is no longer the active/selected one
$jan18 = $wb.Worksheets.Item("jan 18")
$jan18.Activate()
[void]$jan18.Range("B45").Select()

$new.Activate()

Write-Output "done"
